# Adds the 6 newest "Milionária" lottery draws (concursos 314-319) to the
# bottom of the results table on the "+ MILIONÁRIA" sheet, extending the
# used range from A1:I314 to A1:I320, and leaves the new block selected
# the way the author's Excel session ended up (B315:I320, active cell B315).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one drawn row:
#   Concurso, Bola1, Bola2, Bola3, Bola4, Bola5, Bola6, Trevo1, Trevo2
$novosConcursos = @(
    @(314,  7, 12, 14, 17, 21, 41, 4, 5),
    @(315,  5,  7, 11, 12, 20, 29, 2, 4),
    @(316, 10, 11, 12, 18, 36, 48, 3, 4),
    @(317,  8, 23, 28, 40, 49, 50, 1, 2),
    @(318,  3,  6, 22, 23, 38, 49, 4, 5),
    @(319, 11, 18, 26, 29, 30, 37, 1, 3)
)

$primeiraLinha = 315
for ($i = 0; $i -lt $novosConcursos.Count; $i++) {
    $linha = $primeiraLinha + $i
    $valores = $novosConcursos[$i]
    for ($col = 0; $col -lt $valores.Count; $col++) {
        $ws.Cells.Item($linha, $col + 1).Value = $valores[$col]
    }
}

# Match the selection left behind in the workbook after entering the data.
$ws.Range("B315:I320").Select() | Out-Null
